# Insert a new row at position 28 (pushes existing rows 28..76 down to 29..77)
# and populate it with the new weekly price record for "Poroto verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(28).Insert()

$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44629
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112031
$ws.Range("G28").Value = "Poroto verde"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 30000
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 1200
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
